# Apply targeted text replacements to cells F2, H2 and I2 on the active sheet.
# Only the rpc-reply message-id (UUID) values change, plus one flow-id value
# inside H2's "commit" response block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F2: first get-config rpc-reply message-id ---
$f2 = $ws.Range("F2").Value()
$f2 = $f2.Replace("d0480618-c33d-4372-b15e-ae1afb6c733e", "d5857426-5438-4d71-b4f3-dadb5e32ac13")
$ws.Range("F2").Value = $f2

# --- H2: edit-config response message-id, commit response message-id, and flow-id ---
$h2 = $ws.Range("H2").Value()
$h2 = $h2.Replace("844a211c-0605-422d-bfd1-305481269aca", "0a336b10-4bf1-4fd0-a1a1-3f5f979772c7")
$h2 = $h2.Replace("222b944d-cade-4e67-961d-b2a391878e95", "e69b1633-aef2-4e95-bc32-77ddb2a65aa0")
$h2 = $h2.Replace('nc-ext:flow-id="77"', 'nc-ext:flow-id="241"')
$ws.Range("H2").Value = $h2

# --- I2: second get-config rpc-reply message-id ---
$i2 = $ws.Range("I2").Value()
$i2 = $i2.Replace("6ba9f343-5f06-4091-a9eb-83f99db066a5", "8f0c1865-e848-428b-b043-80d562439215")
$ws.Range("I2").Value = $i2
